# "major accuracy check update"
#
# 1) Correct the kit name used throughout column G (NEBNextPoly(A)E7490 -> ...E7490L)
# 2) Re-home the view: scroll back to the top (A1) and move the active
#    selection from column I to column G
# 3) Widen column G so the longer kit name is fully visible

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fix the shared kit-name text in column G (rows 2-49) ---
for ($r = 2; $r -le 49; $r++) {
    $ws.Range("G$r").Value = "NEBNextPoly(A)E7490L"
}

# --- 2. Update the view: scroll to A1 and select G2:G49 ---
$ws.Activate() | Out-Null
$ws.Range("G2:G49").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1

# --- 3. Widen column G to fit the updated text ---
$ws.Columns("G:G").ColumnWidth = 28.86
